$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for several rows per repulled data
$ws.Range("F2").Value = -10
$ws.Range("F3").Value = 4
$ws.Range("F8").Value = 2
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = 2
